$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows (date serial, epidemiological_week, last_available_confirmed,
# last_available_deaths, new_confirmed, new_deaths)
$data = @(
    @(44678, 0, 327052, 6342, 15, 0),
    @(44679, 0, 327065, 6342, 13, 0),
    @(44680, 0, 327071, 6343, 6, 1),
    @(44681, 0, 327076, 6343, 5, 0)
)

$startRow = 32
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]

    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.Value = $values[0]
    $dateCell.NumberFormat = "yyyy\-mm\-dd;@"

    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
}

# Update selection to match the diff (active cell C32)
$ws.Range("C32").Select()
